# Trade #90 closed at 2026-02-17 15:54:11 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.72
$summary.Range("B4").Value = -0.29
$summary.Range("B5").Value = -0.06
$summary.Range("B6").Value = 90
$summary.Range("B7").Value = 31
$summary.Range("B9").Value = 34.44

# --- Strategy Status sheet (row 4 = MarketMaking) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.72
$status.Range("D4").Value = 90
$status.Range("E4").Value = -0.29
$status.Range("F4").Value = -0.28
$status.Range("G4").Value = 34.44

# --- Add new trade row (row 91) to both "All Trades" and "MarketMaking" sheets ---
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(91, 1).Value = 90
    # Dates like "2026-02-17" get smart-parsed into date serials by the
    # COM layer. A leading apostrophe forces text entry (as Excel does
    # for manually-typed values), and resetting the style back to Normal
    # afterwards drops the resulting quote-prefix style so the cell ends
    # up as a plain shared string, matching the rest of the column.
    $ws.Cells.Item(91, 2).Value = "'2026-02-17"
    $ws.Cells.Item(91, 2).Style = "Normal"
    $ws.Cells.Item(91, 3).Value = "15:54:05"
    $ws.Cells.Item(91, 4).Value = "MarketMaking"
    $ws.Cells.Item(91, 5).Value = "UP"
    $ws.Cells.Item(91, 6).Value = 0.68
    $ws.Cells.Item(91, 7).Value = 0.71
    $ws.Cells.Item(91, 8).Value = "CLOSED"
    $ws.Cells.Item(91, 9).Value = 4.4118
    $ws.Cells.Item(91, 10).Value = 0.03
    $ws.Cells.Item(91, 11).Value = 99.72
    $ws.Cells.Item(91, 12).Value = 0
    $ws.Cells.Item(91, 13).Value = 0
    $ws.Cells.Item(91, 14).Value = 0.6
    $ws.Cells.Item(91, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(91, 16).Value = "early_exit"
    $ws.Cells.Item(91, 17).Value = 0.11
}
